$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.597.28'
$ws.Range("E2").Value = '  -3.08%  '
$ws.Range("D3").Value = '2.893.83'
$ws.Range("E3").Value = '  -4.17%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.21'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.11'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.25%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.501'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.39%  '
$ws.Range("D9").Value = '2.890.83'
$ws.Range("E9").Value = '  -4.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.64'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.17%  '
$ws.Range("E11").Value = '  -4.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.446'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000224'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.95%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.81'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.97%  '
$ws.Range("E15").Value = '  +0.51%  '
$ws.Range("D16").Value = '3.373.42'
$ws.Range("E16").Value = '  -4.17%  '
$ws.Range("D17").Value = '60.568.26'
$ws.Range("E17").Value = '  -3.00%  '
$ws.Range("E18").Value = '  -3.29%  '
$ws.Range("D19").Value = '2.894.06'
$ws.Range("E19").Value = '  -4.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '424.44'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.59'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.667'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.27%  '
$ws.Range("E23").Value = '  -5.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.01'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.98'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.91%  '
$ws.Range("E26").Value = '  -1.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.81'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.84%  '
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.22'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.63%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.19'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.94%  '
$ws.Range("E32").Value = '  -3.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.31'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.57%  '
$ws.Range("E34").Value = '  -3.73%  '
$ws.Range("D35").Value = '0.0₃0833'
$ws.Range("E35").Value = '  -2.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.997'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.63'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '49.23'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.73%  '
$ws.Range("E39").Value = '  -2.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.93'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.122'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.68'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.76%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.290'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.56%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.51'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0345'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '371.84'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.12%  '
$ws.Range("D47").Value = '2.646.69'
$ws.Range("E47").Value = '  -3.79%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.76'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.50%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.91'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.96%  '
$ws.Range("E51").Value = '  -1.72%  '
